$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item("Lookup")

# ---------------------------------------------------------------------
# Copy existing header/body formatting onto the new K/L columns so the
# same (reused) style indices are used instead of minting new ones.
# ---------------------------------------------------------------------
$ws2.Range("H1").Copy()
$ws2.Range("K1:L1").PasteSpecial(-4122)

$ws2.Range("H2").Copy()
$ws2.Range("K2:K20").PasteSpecial(-4122)
$ws2.Range("L2:L20").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 1 headers for the new lookup columns
# ---------------------------------------------------------------------
$ws2.Range("K1").Value = "DATA_TYPE"
$ws2.Range("L1").Value = "Input Type"

# ---------------------------------------------------------------------
# Column K: POSTGRES data type (mirrors column H)
# ---------------------------------------------------------------------
$ws2.Range("K2").Value = "INTEGER"
$ws2.Range("K3").Value = "BIGINT"
$ws2.Range("K4").Value = "SMALLINT"
$ws2.Range("K5").Value = "SERIAL"
$ws2.Range("K6").Value = "BIGSERIAL"
$ws2.Range("K7").Value = "BOOLEAN"
$ws2.Range("K8").Value = "CHAR"
$ws2.Range("K9").Value = "VARCHAR"
$ws2.Range("K10").Value = "TEXT"
$ws2.Range("K11").Value = "DATE"
$ws2.Range("K12").Value = "TIMESTAMP"
$ws2.Range("K13").Value = "TIMESTAMPTZ"
$ws2.Range("K14").Value = "REAL"
$ws2.Range("K15").Value = "DOUBLE PRECISION"
$ws2.Range("K16").Value = "NUMERIC"
$ws2.Range("K17").Value = "JSON"
$ws2.Range("K18").Value = "UUID"
$ws2.Range("K19").Value = "BYTEA"
$ws2.Range("K20").Value = "XML"

# ---------------------------------------------------------------------
# Column L: Input Type used to drive form-field generation
# ---------------------------------------------------------------------
$ws2.Range("L2").Value = "Number"
$ws2.Range("L3").Value = "Number"
$ws2.Range("L4").Value = "Number"
$ws2.Range("L5").Value = "Number"
$ws2.Range("L6").Value = "Number"
$ws2.Range("L7").Value = "Radio"
$ws2.Range("L8").Value = "Text"
$ws2.Range("L9").Value = "Text"
$ws2.Range("L10").Value = "TextArea"
$ws2.Range("L11").Value = "Date"
$ws2.Range("L12").Value = "Date"
$ws2.Range("L13").Value = "Date"
$ws2.Range("L14").Value = "Number"
$ws2.Range("L15").Value = "Number"
$ws2.Range("L16").Value = "Number"
$ws2.Range("L17").Value = "TextArea"
$ws2.Range("L18").Value = "Text"
$ws2.Range("L19").Value = "Text"
$ws2.Range("L20").Value = "Text"

# ---------------------------------------------------------------------
# Column N: a little templating helper used to generate a C#-ish enum
# mapping text, with a bold cell + thin left border as a separator
# ---------------------------------------------------------------------
$ws2.Range("N1").Value = '{"$dt", "$it"},'
$ws2.Range("N1").Font.Bold = $true
$ws2.Range("N1").Borders.Item(7).LineStyle = 1

$ws2.Range("N2").Formula = '=SUBSTITUTE(LOWER(SUBSTITUTE($N$1,"$dt", K2)), "$it",L2)'
$ws2.Range("N3:N20").Formula = '=SUBSTITUTE(LOWER(SUBSTITUTE($N$1,"$dt", K3)), "$it",L3)'

# ---------------------------------------------------------------------
# Page setup / view bits: make the Lookup sheet the active tab, select
# the new helper column, and set an explicit print orientation.
# ---------------------------------------------------------------------
$ws2.PageSetup.Orientation = 1

$ws2.Activate()
$ws2.Range("N2:N20").Select()
